$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: tighten the stored precision to 2 decimal places ("custom accuracy").
# Columns that were already at 2 decimals (D, H, W, X, AD) are left untouched.
$roundedValues = @{
    2  = 15.85               # B5
    3  = 11.58                # C5
    5  = 34.45                # E5
    6  = 27.93                # F5
    7  = 12.48                # G5
    9  = 19.2                 # I5
    10 = 8.460000000000001    # J5
    11 = 12.45                # K5
    12 = 13.83                # L5
    13 = 14.55                # M5
    14 = 3.99                 # N5
    15 = 12.41                # O5
    16 = 17.6                 # P5
    17 = 10.55                # Q5
    18 = 0.8                  # R5
    19 = 0.71                 # S5
    20 = 181.53                # T5
    21 = 34.71                # U5
    22 = 11.45                # V5
    25 = 1.96                 # Y5
    26 = 23.63                # Z5
    27 = 10.12                # AA5
    28 = 9                    # AB5
    29 = 10.6                 # AC5
    31 = 0.55                 # AE5
    32 = 44.45                # AF5
    33 = 6.39                 # AG5
    34 = 14.32                # AH5
}

foreach ($col in $roundedValues.Keys) {
    $ws.Cells.Item(5, $col).Value = $roundedValues[$col]
}

# Drop row 6 entirely - the dataset shrinks from A1:AH6 to A1:AH5.
$ws.Rows.Item(6).Delete()
